$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(1)

$origHeight = $shape.Height
$origWidth = $shape.Width
$origTop = $shape.Top
$origLeft = $shape.Left

$lines = @(
    'Title: Identifying the Main Idea of a Passage',
    'Learning Objectives:',
    '- Students will understand the definition of main idea.',
    '- Students will learn how to identify the main idea in a passage.',
    '- Students will be able to use context clues and vocabulary words to infer the main idea.',
    'Starter Activity:',
    '- Ask students to share their understanding of the term "main idea."',
    '- Ask students to share any strategies they use to identify the main idea in a passage.',
    'Worked Examples:',
    '- Display a short passage on the board or overhead.',
    '- Read the passage aloud and underline key words and phrases to identify the main idea.',
    '- Have students work in pairs to read a different passage and identify the main idea together.',
    '- Discuss the main ideas identified by each pair as a class.',
    'Practice:',
    '- Distribute a worksheet to students with several short passages.',
    '- Instruct students to read each passage and circle the sentence that represents the main idea.',
    '- Review their answers as a class to identify any common themes or difficulties.',
    'Exit Ticket:',
    '- Students will write the main idea of a brief passage on a sticky note and stick it to the board as they exit the classroom.',
    'Materials needed:',
    '- Powerpoint presentation with the title and learning objectives included.',
    '- Short passages for worked examples and practice worksheet.',
    '- Sticky notes for exit ticket activity.'
)

$newText = [string]::Join([char]13, $lines)

$shape.TextFrame.TextRange.Text = $newText

$shape.Height = $origHeight
$shape.Width = $origWidth
$shape.Top = $origTop
$shape.Left = $origLeft
